$d = $word.ActiveDocument

$replacements = @(
    @{old="177×9="; new="477×9="},
    @{old="746×5="; new="472×7="},
    @{old="388×3="; new="591×7="},
    @{old="928×2="; new="987×8="},
    @{old="547×3="; new="652×4="},
    @{old="922×2="; new="511×7="},
    @{old="573×8="; new="354×7="},
    @{old="453×4="; new="252×7="},
    @{old="963×8="; new="959×3="},
    @{old="373×2="; new="148×8="},
    @{old="219×4="; new="641×8="},
    @{old="413×6="; new="232×8="},
    @{old="171×2="; new="322×4="},
    @{old="629×2="; new="689×2="},
    @{old="178×8="; new="498×9="},
    @{old="618×4="; new="889×5="},
    @{old="831×9="; new="212×6="},
    @{old="762×4="; new="626×7="},
    @{old="956×9="; new="222×4="},
    @{old="734×9="; new="114×7="},
    @{old="247×2="; new="820×3="},
    @{old="314×2="; new="596×6="},
    @{old="511×2="; new="694×6="},
    @{old="651×7="; new="284×8="},
    @{old="575×8="; new="637×2="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
